$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$origStyle = $ws.Range('D2').Style
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.816.68'
$ws.Range('D2').Style = $origStyle
$ws.Range('E2').Value = '  +1.81%  '

# Row 3
$origStyle = $ws.Range('D3').Style
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.284.23'
$ws.Range('D3').Style = $origStyle
$ws.Range('E3').Value = '  +0.91%  '

# Row 4
$ws.Range('E4').Value = '  +0.06%  '

# Row 5
$origStyle = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '584.27'
$ws.Range('D5').Style = $origStyle
$ws.Range('E5').Value = '  +1.08%  '

# Row 6
$origStyle = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '182.30'
$ws.Range('D6').Style = $origStyle
$ws.Range('E6').Value = '  +0.33%  '

# Row 7
$ws.Range('E7').Value = '  +0.03%  '

# Row 8
$ws.Range('E8').Value = '  +0.34%  '

# Row 9
$ws.Range('E9').Value = '  +2.24%  '

# Row 10
$origStyle = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.68'
$ws.Range('D10').Style = $origStyle
$ws.Range('E10').Value = '  -1.66%  '

# Row 11
$origStyle = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.428'
$ws.Range('D11').Style = $origStyle
$ws.Range('E11').Value = '  +3.19%  '

# Row 12
$origStyle = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '3.851.21'
$ws.Range('D12').Style = $origStyle
$ws.Range('E12').Value = '  +0.90%  '

# Row 13
$ws.Range('E13').Value = '  +0.03%  '

# Row 14
$origStyle = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.89'
$ws.Range('D14').Style = $origStyle
$ws.Range('E14').Value = '  +0.37%  '

# Row 15
$origStyle = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '68.845.99'
$ws.Range('D15').Style = $origStyle
$ws.Range('E15').Value = '  +1.89%  '

# Row 16
$ws.Range('E16').Value = '  +2.57%  '

# Row 17
$origStyle = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.240.03'
$ws.Range('D17').Style = $origStyle
$ws.Range('E17').Value = '  -0.54%  '

# Row 18
$ws.Range('E18').Value = '  +0.16%  '

# Row 19
$origStyle = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.62'
$ws.Range('D19').Style = $origStyle
$ws.Range('E19').Value = '  +0.51%  '

# Row 20
$origStyle = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '395.46'
$ws.Range('D20').Style = $origStyle
$ws.Range('E20').Value = '  +4.93%  '

# Row 21
$origStyle = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.72'
$ws.Range('D21').Style = $origStyle
$ws.Range('E21').Value = '  +1.09%  '

# Row 22
$origStyle = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '72.10'
$ws.Range('D22').Style = $origStyle

# Row 23
$ws.Range('E23').Value = '  +0.03%  '

# Row 24
$origStyle = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.517'
$ws.Range('D24').Style = $origStyle
$ws.Range('E24').Value = '  +0.87%  '

# Row 25
$origStyle = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000120'
$ws.Range('D25').Style = $origStyle
$ws.Range('E25').Value = '  +0.47%  '

# Row 26
$ws.Range('E26').Value = '  +4.98%  '

# Row 27
$origStyle = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.73'
$ws.Range('D27').Style = $origStyle
$ws.Range('E27').Value = '  +0.89%  '

# Row 28
$ws.Range('E28').Value = '  -0.50%  '

# Row 29
$origStyle = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.98'
$ws.Range('D29').Style = $origStyle
$ws.Range('E29').Value = '  +0.04%  '

# Row 30
$origStyle = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.72'
$ws.Range('D30').Style = $origStyle
$ws.Range('E30').Value = '  -1.06%  '

# Row 31
$origStyle = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '23.04'
$ws.Range('D31').Style = $origStyle
$ws.Range('E31').Value = '  +1.45%  '

# Row 32
$ws.Range('B32').Value = 'Aptos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$origStyle = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.16'
$ws.Range('D32').Style = $origStyle
$ws.Range('E32').Value = '  +3.13%  '

# Row 33
$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$origStyle = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.30'
$ws.Range('D33').Style = $origStyle
$ws.Range('E33').Value = '  +1.86%  '

# Row 34
$ws.Range('E34').Value = '  +0.06%  '

# Row 35
$ws.Range('B35').Value = 'Monero'
$ws.Range('C35').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$origStyle = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '164.61'
$ws.Range('D35').Style = $origStyle
$ws.Range('E35').Value = '  +0.53%  '

# Row 36
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$origStyle = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.51'
$ws.Range('D36').Style = $origStyle
$ws.Range('E36').Value = '  +0.73%  '

# Row 37
$ws.Range('E37').Value = '  +3.79%  '

# Row 38
$ws.Range('E38').Value = '  -2.23%  '

# Row 39
$origStyle = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '26.44'
$ws.Range('D39').Style = $origStyle
$ws.Range('E39').Value = '  -1.54%  '

# Row 40
$ws.Range('E40').Value = '  +0.76%  '

# Row 41
$origStyle = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.59'
$ws.Range('D41').Style = $origStyle
$ws.Range('E41').Value = '  -2.76%  '

# Row 42
$ws.Range('E42').Value = '  -2.59%  '

# Row 43
$ws.Range('E43').Value = '  +1.32%  '

# Row 44
$ws.Range('B44').Value = 'Hedera'
$ws.Range('C44').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$origStyle = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0690'
$ws.Range('D44').Style = $origStyle
$ws.Range('E44').Value = '  +1.01%  '

# Row 45
$ws.Range('B45').Value = 'Bittensor'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$origStyle = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '345.97'
$ws.Range('D45').Style = $origStyle
$ws.Range('E45').Value = '  -4.91%  '

# Row 46
$origStyle = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.627.90'
$ws.Range('D46').Style = $origStyle
$ws.Range('E46').Value = '  -4.07%  '

# Row 47
$origStyle = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '24.71'
$ws.Range('D47').Style = $origStyle
$ws.Range('E47').Value = '  -2.81%  '

# Row 48
$origStyle = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0283'
$ws.Range('D48').Style = $origStyle
$ws.Range('E48').Value = '  +1.16%  '

# Row 49
$origStyle = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '31.96'
$ws.Range('D49').Style = $origStyle
$ws.Range('E49').Value = '  +1.67%  '

# Row 50
$ws.Range('E50').Value = '  +2.87%  '

# Row 51
$ws.Range('E51').Value = '  +0.04%  '
